$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire column A (the GENE index numbers), shifting
# columns B:F left to become A:E.
$ws.Range("A:A").Delete()
